# Update the "UnitMass" column (column C) values in Sheet1 for the
# LoadingPC1 loadings table (both the "+ loading" block, rows 2-21,
# and the "- loading" block, rows 23-42), matching the new
# hyperparameter-optimization sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 39
    3  = 27
    4  = 28
    5  = 40
    6  = 24
    7  = 25
    8  = 23
    9  = 1
    10 = 7
    11 = 46
    12 = 26
    13 = 56
    14 = 57
    15 = 45
    16 = 20
    17 = 117
    18 = 16
    19 = 14
    20 = 12
    21 = 6
    23 = 43
    24 = 70
    25 = 69
    26 = 55
    27 = 87
    28 = 71
    29 = 86
    30 = 73
    31 = 72
    32 = 59
    33 = 60
    34 = 58
    35 = 85
    36 = 84
    37 = 81
    38 = 97
    39 = 19
    40 = 31
    41 = 53
    42 = 42
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
